{"js": "// Replace the 25 two-digit multiplication expressions in the table cells.\n// Each (old, new) pair is applied with an exact, case-sensitive search\n// against the document body; since every old expression is unique in the\n// document, a single search+replace per pair is unambiguous.\nconst replacements = [\n  [\"63\u00d786=\", \"17\u00d767=\"],\n  [\"76\u00d768=\", \"41\u00d755=\"],\n  [\"17\u00d728=\", \"29\u00d756=\"],\n  [\"41\u00d725=\", \"58\u00d754=\"],\n  [\"68\u00d729=\", \"96\u00d737=\"],\n  [\"50\u00d718=\", \"34\u00d734=\"],\n  [\"68\u00d728=\", \"50\u00d739=\"],\n  [\"81\u00d716=\", \"75\u00d788=\"],\n  [\"23\u00d797=\", \"18\u00d770=\"],\n  [\"46\u00d762=\", \"47\u00d778=\"],\n  [\"68\u00d796=\", \"82\u00d796=\"],\n  [\"82\u00d798=\", \"57\u00d747=\"],\n  [\"26\u00d731=\", \"58\u00d726=\"],\n  [\"26\u00d781=\", \"87\u00d770=\"],\n  [\"23\u00d711=\", \"36\u00d758=\"],\n  [\"59\u00d762=\", \"19\u00d762=\"],\n  [\"63\u00d796=\", \"70\u00d799=\"],\n  [\"71\u00d748=\", \"45\u00d735=\"],\n  [\"99\u00d771=\", \"99\u00d759=\"],\n  [\"15\u00d743=\", \"79\u00d725=\"],\n  [\"81\u00d754=\", \"23\u00d759=\"],\n  [\"48\u00d762=\", \"12\u00d791=\"],\n  [\"80\u00d713=\", \"58\u00d775=\"],\n  [\"47\u00d736=\", \"86\u00d777=\"],\n  [\"31\u00d765=\", \"97\u00d792=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication expressions in the table cells.\n# Each (old, new) pair is applied with Word's Find/Replace against the\n# whole document content; every old expression is unique in the document,\n# so MatchCase + ReplaceAll is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"63\u00d786=\", \"17\u00d767=\"),\n    @(\"76\u00d768=\", \"41\u00d755=\"),\n    @(\"17\u00d728=\", \"29\u00d756=\"),\n    @(\"41\u00d725=\", \"58\u00d754=\"),\n    @(\"68\u00d729=\", \"96\u00d737=\"),\n    @(\"50\u00d718=\", \"34\u00d734=\"),\n    @(\"68\u00d728=\", \"50\u00d739=\"),\n    @(\"81\u00d716=\", \"75\u00d788=\"),\n    @(\"23\u00d797=\", \"18\u00d770=\"),\n    @(\"46\u00d762=\", \"47\u00d778=\"),\n    @(\"68\u00d796=\", \"82\u00d796=\"),\n    @(\"82\u00d798=\", \"57\u00d747=\"),\n    @(\"26\u00d731=\", \"58\u00d726=\"),\n    @(\"26\u00d781=\", \"87\u00d770=\"),\n    @(\"23\u00d711=\", \"36\u00d758=\"),\n    @(\"59\u00d762=\", \"19\u00d762=\"),\n    @(\"63\u00d796=\", \"70\u00d799=\"),\n    @(\"71\u00d748=\", \"45\u00d735=\"),\n    @(\"99\u00d771=\", \"99\u00d759=\"),\n    @(\"15\u00d743=\", \"79\u00d725=\"),\n    @(\"81\u00d754=\", \"23\u00d759=\"),\n    @(\"48\u00d762=\", \"12\u00d791=\"),\n    @(\"80\u00d713=\", \"58\u00d775=\"),\n    @(\"47\u00d736=\", \"86\u00d777=\"),\n    @(\"31\u00d765=\", \"97\u00d792=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        [ref]$oldText, $true, $true, $false, $false, $false,\n        $true, 1, $false, [ref]$newText, 2\n    ) | Out-Null\n}\n"}
